$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Result")

$ws.Range("C2").Value = "2015-09-17 11:38:23"
$ws.Range("D2").Value = "2015-09-17 11:38:25"
